$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1340.6842
$ws.Range("I2").Value = 205.09091
$ws.Range("K2").Value = 205.09091
$ws.Range("M2").Value = -92.09091000000001
$ws.Range("H19").Value = 2607.7778
$ws.Range("I19").Value = 1495
$ws.Range("K19").Value = 1495
$ws.Range("M19").Value = -1320
$ws.Range("H62").Value = 7944227.5
$ws.Range("I62").Value = 8936381
$ws.Range("J62").Value = 6999.5
$ws.Range("K62").Value = 8936381
$ws.Range("L62").Value = 6999.5
$ws.Range("M62").Value = -8935757
$ws.Range("N62").Value = -8247.5
$ws.Range("H65").Value = 7944227.5
$ws.Range("I65").Value = 8936381
$ws.Range("J65").Value = 6999.5
$ws.Range("K65").Value = 44681905
$ws.Range("L65").Value = 34997.5
$ws.Range("M65").Value = -44678785
$ws.Range("N65").Value = -41237.5
$ws.Range("H98").Value = 2416.7827
$ws.Range("I98").Value = 2512.5454
$ws.Range("J98").Value = 310
$ws.Range("K98").Value = 2512.5454
$ws.Range("L98").Value = 310
$ws.Range("M98").Value = -1014.5454
$ws.Range("N98").Value = -3306
$ws.Range("H122").Value = 2416.7827
$ws.Range("I122").Value = 2512.5454
$ws.Range("J122").Value = 310
$ws.Range("K122").Value = 7537.6362
$ws.Range("L122").Value = 930
$ws.Range("M122").Value = -5087.6362
$ws.Range("N122").Value = -5830

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3453.8928
$ws.Range("I2").Value = 1845.3334
$ws.Range("J2").Value = 5309.923
$ws.Range("K2").Value = 1845.3334
$ws.Range("L2").Value = 5309.923
$ws.Range("M2").Value = -1732.3334
$ws.Range("N2").Value = -5535.923
$ws.Range("H32").Value = 2272118.5
$ws.Range("I32").Value = 3597.3684
$ws.Range("K32").Value = 3597.3684
$ws.Range("M32").Value = -3310.3684
$ws.Range("H76").Value = 18800
$ws.Range("J76").Value = 18800
$ws.Range("L76").Value = 18800
$ws.Range("N76").Value = -19476
$ws.Range("H79").Value = 18800
$ws.Range("J79").Value = 18800
$ws.Range("L79").Value = 18800
$ws.Range("N79").Value = -21140
$ws.Range("H109").Value = 22000
$ws.Range("J109").Value = 22000
$ws.Range("L109").Value = 22000
$ws.Range("N109").Value = -24774
$ws.Range("H112").Value = 123998
$ws.Range("J112").Value = 123998
$ws.Range("L112").Value = 123998
$ws.Range("N112").Value = -126952
$ws.Range("H116").Value = 3453.8928
$ws.Range("I116").Value = 1845.3334
$ws.Range("J116").Value = 5309.923
$ws.Range("K116").Value = 1845.3334
$ws.Range("L116").Value = 5309.923
$ws.Range("M116").Value = 448.6666
$ws.Range("N116").Value = -9897.922999999999
$ws.Range("H132").Value = 1246737.9
$ws.Range("I132").Value = 1424656.1
$ws.Range("J132").Value = 268187.5
$ws.Range("K132").Value = 4273968.300000001
$ws.Range("L132").Value = 804562.5
$ws.Range("M132").Value = -4271438.300000001
$ws.Range("N132").Value = -809622.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3453.8928
$ws.Range("I3").Value = 1845.3334
$ws.Range("J3").Value = 5309.923
$ws.Range("K3").Value = 1845.3334
$ws.Range("L3").Value = 5309.923
$ws.Range("M3").Value = -1731.3334
$ws.Range("N3").Value = -5537.923
$ws.Range("H99").Value = 7231.041
$ws.Range("J99").Value = 7065.871
$ws.Range("L99").Value = 7065.871
$ws.Range("N99").Value = -10061.871

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 6372.0625
$ws.Range("I7").Value = 5721.75
$ws.Range("K7").Value = 5721.75
$ws.Range("M7").Value = -5608.75
$ws.Range("H22").Value = 1786498.8
$ws.Range("I22").Value = 4464710
$ws.Range("K22").Value = 4464710
$ws.Range("M22").Value = -4464360
$ws.Range("H31").Value = 5603.2964
$ws.Range("I31").Value = 2816.1667
$ws.Range("J31").Value = 6399.619
$ws.Range("K31").Value = 2816.1667
$ws.Range("L31").Value = 6399.619
$ws.Range("M31").Value = -2521.1667
$ws.Range("N31").Value = -6989.619
$ws.Range("H34").Value = 5603.2964
$ws.Range("I34").Value = 2816.1667
$ws.Range("J34").Value = 6399.619
$ws.Range("K34").Value = 2816.1667
$ws.Range("L34").Value = 6399.619
$ws.Range("M34").Value = -2614.1667
$ws.Range("N34").Value = -6803.619
$ws.Range("H99").Value = 18522528
$ws.Range("I99").Value = 27781330
$ws.Range("K99").Value = 27781330
$ws.Range("M99").Value = -27779832
$ws.Range("H107").Value = 978.5
$ws.Range("I107").Value = 465.8
$ws.Range("K107").Value = 465.8
$ws.Range("M107").Value = 1454.2
$ws.Range("H126").Value = 18522528
$ws.Range("I126").Value = 27781330
$ws.Range("K126").Value = 83343990
$ws.Range("M126").Value = -83341520

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 1046
$ws.Range("I16").Value = 990
$ws.Range("J16").Value = 1102
$ws.Range("K16").Value = 2970
$ws.Range("L16").Value = 3306
$ws.Range("M16").Value = -2797
$ws.Range("N16").Value = -3652
$ws.Range("H70").Value = 25000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 25000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H131").Value = 32101884
$ws.Range("I131").Value = 44447732
$ws.Range("K131").Value = 133343196
$ws.Range("M131").Value = -133338156

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 39000
$ws.Range("J49").Value = 39000
$ws.Range("L49").Value = 39000
$ws.Range("N49").Value = -39368

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H82").Value = 5937.0835
$ws.Range("I82").Value = 1624.375
$ws.Range("J82").Value = 14562.5
$ws.Range("K82").Value = 1624.375
$ws.Range("L82").Value = 14562.5
$ws.Range("M82").Value = -1263.375
$ws.Range("N82").Value = -15284.5
$ws.Range("H85").Value = 5937.0835
$ws.Range("I85").Value = 1624.375
$ws.Range("J85").Value = 14562.5
$ws.Range("K85").Value = 1624.375
$ws.Range("L85").Value = 14562.5
$ws.Range("M85").Value = -376.375
$ws.Range("N85").Value = -17058.5
$ws.Range("H100").Value = 3771.2144
$ws.Range("I100").Value = 7949.25
$ws.Range("K100").Value = 7949.25
$ws.Range("M100").Value = -7408.25
$ws.Range("H132").Value = 2884
$ws.Range("I132").Value = 2661.2068
$ws.Range("K132").Value = 7983.6204
$ws.Range("M132").Value = -5453.6204

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 26639.715
$ws.Range("I69").Value = 29989.5
$ws.Range("J69").Value = 25299.8
$ws.Range("K69").Value = 29989.5
$ws.Range("L69").Value = 25299.8
$ws.Range("M69").Value = -29240.5
$ws.Range("N69").Value = -26797.8
$ws.Range("H72").Value = 26639.715
$ws.Range("I72").Value = 29989.5
$ws.Range("J72").Value = 25299.8
$ws.Range("K72").Value = 89968.5
$ws.Range("L72").Value = 75899.39999999999
$ws.Range("M72").Value = -86224.5
$ws.Range("N72").Value = -83387.39999999999
$ws.Range("H126").Value = 8633.038
$ws.Range("I126").Value = 7094.6665
$ws.Range("K126").Value = 21283.9995
$ws.Range("M126").Value = -18813.9995
$ws.Range("H132").Value = 4697.143
$ws.Range("I132").Value = 4151.064
$ws.Range("J132").Value = 7548.8887
$ws.Range("K132").Value = 12453.192
$ws.Range("L132").Value = 22646.6661
$ws.Range("M132").Value = -9923.192000000001
$ws.Range("N132").Value = -27706.6661
